$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet ("ODI Batting Extra").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$data = @(
    @("3305", "0", "20.00%"),
    @("3321", "", ""),
    @("3322", "0", "10.00%"),
    @("3323", "", ""),
    @("3325", "1", "10.00%"),
    @("3331", "", ""),
    @("3335", "0", "10.00%"),
    @("3337", "0", ""),
    @("3341", "", ""),
    @("3383", "0", "40.00%"),
    @("3418", "", ""),
    @("3440", "0", "30.00%"),
    @("3442", "", ""),
    @("3444", "0", "10.00%"),
    @("3459", "", ""),
    @("3461", "0", "20.00%"),
    @("3463", "0", ""),
    @("3468", "0", "10.00%"),
    @("3501", "1", "10.00%"),
    @("3503", "0", "")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c]
    }
}
